# Generate Report for Handoff
#
# The handoff/handback process re-ran for the rows whose source files are
# still cycling through the pipeline (everything except the two rows that
# are already fully "Handed back: in sync with en-US" and the one row that
# is still "In Translation"). That produced new timestamps for:
#   - Overview!D  (Latest Handoff Date)
#   - zh-cn!E     (Latest Handoff Datetime)
#   - de-de!E     (Latest Handoff Datetime)
# on rows 4, 6, 7, 8, 9 and 10.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$rows = @(4, 6, 7, 8, 9, 10)

foreach ($r in $rows) {
    $overview.Range("D$r").Value = "2016-03-21 05:26:15"
    $zhcn.Range("E$r").Value     = "2016-03-21 05:26:06"
    $dede.Range("E$r").Value     = "2016-03-21 05:26:15"
}
